$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking text (e.g. "1.000", "20.60",
# "0.00001066") that must stay exact text - Excel's COM layer auto-coerces
# plain numeric strings into real numbers (losing trailing zeros / using
# scientific notation), so every Price assignment is forced to text with a
# leading apostrophe, matching how Excel itself preserves literal text entry.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'28.191.43"
$ws.Range("E2").Value = "  -1.42%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.805.44"
$ws.Range("E3").Value = "  +0.51%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'316.84"
$ws.Range("E5").Value = "  +1.10%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.11%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.5409"
$ws.Range("E7").Value = "  +1.74%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3787"
$ws.Range("E8").Value = "  +0.54%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.07485"
$ws.Range("E9").Value = "  -0.55%  "

# Row 10 - OKB
$ws.Range("D10").Value = "'42.15"
$ws.Range("E10").Value = "  -1.03%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'1.098"
$ws.Range("E11").Value = "  -1.93%  "

# Row 12 - BinanceUSD
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.09%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'6.207"
$ws.Range("E13").Value = "  +0.09%  "

# Row 14 - Solana
$ws.Range("D14").Value = "'20.60"
$ws.Range("E14").Value = "  -2.57%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'7.374"
$ws.Range("E15").Value = "  -1.42%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "'1.805.21"
$ws.Range("E16").Value = "  +0.57%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'90.11"
$ws.Range("E17").Value = "  -0.28%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.00001066"
$ws.Range("E18").Value = "  -0.24%  "

# Row 19 - TRON
$ws.Range("D19").Value = "'0.06512"

# Row 20 - Avalanche
$ws.Range("E20").Value = "  +0.89%  "

# Row 21 - Dai
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.943"
$ws.Range("E22").Value = "  +0.43%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "'28.219.07"
$ws.Range("E23").Value = "  -1.40%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  +0.22%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.090"
$ws.Range("E25").Value = "  -0.41%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'156.03"
$ws.Range("E26").Value = "  -3.03%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'20.54"
$ws.Range("E27").Value = "  -0.13%  "

# Row 28 - WrappedliquidstakedEther2.0
$ws.Range("D28").Value = "'2.011.93"
$ws.Range("E28").Value = "  +0.51%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.343"
$ws.Range("E29").Value = "  -0.89%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'122.14"
$ws.Range("E30").Value = "  -1.24%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.1119"
$ws.Range("E31").Value = "  +9.11%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'1.127"
$ws.Range("E32").Value = "  +0.53%  "

# Row 33 - was Filecoin, now HuobiToken (rows 33/34 swap their row-relative content)
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.674"
$ws.Range("E33").Value = "  +0.13%  "

# Row 34 - was HuobiToken, now Filecoin
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.599"
$ws.Range("E34").Value = "  -2.04%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "'0.06957"
$ws.Range("E35").Value = "  +6.38%  "

# Row 36 - Algorand
$ws.Range("D36").Value = "'0.2231"
$ws.Range("E36").Value = "  -3.03%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.02307"
$ws.Range("E37").Value = "  -0.69%  "

# Row 38 - InternetComputer(DFINITY)
$ws.Range("D38").Value = "'5.110"
$ws.Range("E38").Value = "  +0.89%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "'8.488"
$ws.Range("E39").Value = "  -3.73%  "

# Row 40 - was Aptos, now TheSandbox (rows 40/41 swap their row-relative content)
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6192"
$ws.Range("E40").Value = "  -1.74%  "

# Row 41 - was TheSandbox, now Aptos
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'11.16"
$ws.Range("E41").Value = "  -2.80%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "'1.176"
$ws.Range("E42").Value = "  -2.50%  "

# Row 43 - WEMIXTOKEN
$ws.Range("D43").Value = "'1.423"
$ws.Range("E43").Value = "  +2.01%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "'13.51"
$ws.Range("E44").Value = "  -0.41%  "

# Row 45 - PancakeSwap
$ws.Range("D45").Value = "'3.687"
$ws.Range("E45").Value = "  +0.57%  "

# Row 46 - Decentraland
$ws.Range("D46").Value = "'0.5780"
$ws.Range("E46").Value = "  -2.41%  "

# Row 47 - Quant
$ws.Range("D47").Value = "'125.44"
$ws.Range("E47").Value = "  -0.64%  "

# Row 48 - EOS
$ws.Range("D48").Value = "'1.191"
$ws.Range("E48").Value = "  +1.83%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "'1.931"
$ws.Range("E49").Value = "  -2.42%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "'0.06824"
$ws.Range("E50").Value = "  -1.49%  "

# Row 51 - Aave
$ws.Range("D51").Value = "'72.03"
$ws.Range("E51").Value = "  -1.59%  "
